$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new timesheet entry as row 28 (XSLConstructor v2. outputTemplate jako embedded resource)
$ws.Range("A28").Value = 44028
$ws.Range("B28").Value = 2
$ws.Range("C28").Value = "XSLConstructor. Embedded resources."

# Reflect where the user's selection ended up after entering the new row
$ws.Range("C29").Select()

$wb.Save()
